$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("E18").Value = "ELM-2NA-Máquinas Elétricas"
$ws.Range("F18").Value = "['ELM-2NA-Instalções Elétricas', -]"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("E19").Value = "ELM-2NA-Máquinas Elétricas"
$ws.Range("F19").Value = "['ELM-2NA-Instalções Elétricas', -]"

# Row 20
$ws.Range("E20").Value = "-"
